$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '24.420.50'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.41%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.682.90'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9996'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '316.33'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9993'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3881'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.21%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4002'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.93%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.481'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.55%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.9998'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.32%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '52.27'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.15%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08732'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '25.74'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +10.06%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.477'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.17%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.971'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.670.99'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '97.66'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.60%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.07196'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.57%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '19.62'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.241'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.69%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9996'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('E23').Value = '  -2.35%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '24.405.32'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.41%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.013'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -7.28%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.338'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.73%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.48'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.38%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '167.55'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.85%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.638'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +11.35%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.352'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.48%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '137.75'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.00%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.854.53'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.66%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.08741'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '7.347'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.84%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.042'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.96%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02974'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +8.76%  '
$ws.Range('B37').Value = 'WEMIXTOKEN'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.969'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2741'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '10.77'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.32%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.09127'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.72%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '14.03'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.35%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.7945'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.55%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.470'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '17.29'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +8.93%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.7181'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.577'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.258'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.396'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +6.65%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.9994'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.34%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '139.29'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.16%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.08040'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.75%  '
